# "format date in excel"
#
# The "Modified" column of Table1 holds serial date-time numbers
# (e.g. 43350.0950182755) that are currently displayed with the
# default "General" number format. Format them as dates (dd/MM/yyyy)
# so they render as readable dates instead of raw serial numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefer going through the table/ListObject so the format lands on
# exactly the "Modified" column's data cells, however the sheet is
# shaped.
$lo = $ws.ListObjects.Item("Table1")
$col = $lo.ListColumns.Item("Modified")

$col.DataBodyRange.NumberFormat = "dd/MM/yyyy"
